# Updates the "Price" (D) and "Volume(1h)" (E) columns of the crypto
# tracker sheet with freshly scraped values, matching the automated
# "Updated cryptos list ... with GitHub Actions" commit.
#
# Note: several Price values look like plain numbers (e.g. "73.05",
# "166.90") but must remain text so trailing zeros / exact formatting
# are preserved. Those are written with a leading apostrophe
# (PowerShell: "`'") which Excel interprets as "force text", matching
# how the sheet already stores all other cells as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.952.38"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "2.214.27"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "`'240.53"
$ws.Range("E5").Value = "  -2.67%  "
$ws.Range("D6").Value = "`'0.624"
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("D7").Value = "`'73.21"
$ws.Range("E7").Value = "  -1.78%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").Value = "`'0.606"
$ws.Range("D10").Value = "`'42.46"
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("D11").Value = "`'0.0951"
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").Value = "`'7.07"
$ws.Range("E12").Value = "  -1.85%  "
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("D14").Value = "2.548.42"
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("D15").Value = "`'14.26"
$ws.Range("D16").Value = "`'0.836"
$ws.Range("E16").Value = "  -1.97%  "
$ws.Range("D17").Value = "2.215.86"
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").Value = "41.870.48"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").Value = "`'0.0000106"
$ws.Range("E19").Value = "  +7.62%  "
$ws.Range("D20").Value = "`'73.05"
$ws.Range("D21").Value = "`'6.15"
$ws.Range("E22").Value = "  +20.92%  "
$ws.Range("D23").Value = "`'229.18"
$ws.Range("E23").Value = "  -1.08%  "
$ws.Range("E24").Value = "  -6.39%  "
$ws.Range("D25").Value = "`'11.72"
$ws.Range("E25").Value = "  +2.64%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").Value = "`'2.26"
$ws.Range("E28").Value = "  -2.05%  "
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("D30").Value = "`'166.90"
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("D31").Value = "`'20.43"
$ws.Range("E31").Value = "  -1.29%  "
$ws.Range("D32").Value = "`'5.58"
$ws.Range("E32").Value = "  +6.32%  "
$ws.Range("D33").Value = "`'0.0794"
$ws.Range("E33").Value = "  -4.16%  "
$ws.Range("D34").Value = "`'0.124"
$ws.Range("E34").Value = "  -0.63%  "
$ws.Range("D35").Value = "`'29.30"
$ws.Range("E35").Value = "  -6.28%  "
$ws.Range("E36").Value = "  -10.91%  "
$ws.Range("D37").Value = "`'4.26"
$ws.Range("E37").Value = "  -5.24%  "
$ws.Range("D38").Value = "`'0.0298"
$ws.Range("E38").Value = "  -5.29%  "
$ws.Range("E39").Value = "  -1.04%  "
$ws.Range("D40").Value = "`'65.77"
$ws.Range("E40").Value = "  +5.91%  "
$ws.Range("E41").Value = "  -3.08%  "
$ws.Range("D42").Value = "`'5.61"
$ws.Range("E42").Value = "  -2.91%  "
$ws.Range("D43").Value = "`'0.197"
$ws.Range("E43").Value = "  -3.58%  "
$ws.Range("D44").Value = "`'8.69"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").Value = "`'103.99"
$ws.Range("E45").Value = "  -2.92%  "
$ws.Range("E46").Value = "  -3.03%  "
$ws.Range("E47").Value = "  +2.68%  "
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").Value = "`'2.70"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").Value = "2.423.32"
$ws.Range("E51").Value = "  -1.48%  "
